# Add files via upload
# Marks US01/US02 (Sprint1) as Completed with actual size/time/completion
# date, renames US02 to reflect paired programming and updates its
# owner pairing, and syncs the Backlog status for both stories.

$wb = $excel.ActiveWorkbook

# --- Sprint1 sheet -------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")

# Row 5 -> US01 "Date before current dates"
# (workbook uses the 1904 date system; 42646 == 10/4/2020 in that system,
# matching the completion date already used on rows 3 & 4)
$sprint1.Range("D5").Value = "Completed"
$sprint1.Range("G5").Value = 15
$sprint1.Range("H5").Value = 60
$sprint1.Range("I5").Value = 42646

# Row 6 -> US02, renamed + paired programming owners
$sprint1.Range("B6").Value = "Birth before marriage(Paired Programming)"
$sprint1.Range("C6").Value = "AA/DA"
$sprint1.Range("D6").Value = "Completed"
$sprint1.Range("G6").Value = 15
$sprint1.Range("H6").Value = 60
$sprint1.Range("I6").Value = 42646

# --- Backlog sheet ---------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("E10").Value = "Completed"
$backlog.Range("E11").Value = "Completed"

# --- Active sheet moved to Stories prior to save --------------------------
$wb.Worksheets.Item("Stories").Activate()
